# Mass Interview uncheck for Only Available Interviewers
# Fills in newly-run interview/sprint rows that were previously blank
# placeholder rows in the "AMSIN" and "BETA" sheets.

$wb = $excel.ActiveWorkbook

# ---- AMSIN sheet: row 24 gets its timestamp nudged, and two new rows
#      (25, 26) of sprint data are appended after it. ----
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$wsAmsin.Cells.Item(24, 2).Value = 44476.38354539352

$wsAmsin.Cells.Item(25, 1).Value = "2021-10-26"
$wsAmsin.Cells.Item(25, 2).Value = 44495.66112016204
$wsAmsin.Cells.Item(25, 3).Value = "152_fstcycle"
$wsAmsin.Cells.Item(25, 4).Value = 89
$wsAmsin.Cells.Item(25, 5).Value = 88
$wsAmsin.Cells.Item(25, 6).Value = 1
$wsAmsin.Cells.Item(25, 7).Value = 2.95

$wsAmsin.Cells.Item(26, 1).Value = "2021-10-28"
$wsAmsin.Cells.Item(26, 2).Value = 44497.40656510417
$wsAmsin.Cells.Item(26, 3).Value = "152_fnlrgrsn"
$wsAmsin.Cells.Item(26, 4).Value = 89
$wsAmsin.Cells.Item(26, 5).Value = 89
$wsAmsin.Cells.Item(26, 6).Value = 0
$wsAmsin.Cells.Item(26, 7).Value = 2.49

# ---- BETA sheet: row 13 was a blank placeholder row; fill it with the
#      new interview run data. ----
$wsBeta = $wb.Worksheets.Item("BETA")

$wsBeta.Cells.Item(13, 1).Value = "2021-10-28"
$wsBeta.Cells.Item(13, 2).Value = 44497.70566909626
$wsBeta.Cells.Item(13, 3).Value = "152_betachgs"
$wsBeta.Cells.Item(13, 4).Value = 89
$wsBeta.Cells.Item(13, 5).Value = 89
$wsBeta.Cells.Item(13, 6).Value = 0
$wsBeta.Cells.Item(13, 7).Value = 4.23
